$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 56 (shifts old rows 56-61 down to 57-62) for the
# new GUADALUPE port entry that was inserted in the middle of the table.
$ws.Rows.Item(56).Insert()

# Populate the port-name (text) cells in alphabetical order so that the
# underlying shared-string table is built up in that same order, matching
# how the source workbook ended up with its new unique strings appended.
$ws.Range("A63").Value = "CHINOOK"
$ws.Range("A64").Value = "FIELDS LANDING"
$ws.Range("A56").Value = "GUADALUPE"
$ws.Range("A65").Value = "MARCONI"
$ws.Range("A66").Value = "PORT ORFORD"
$ws.Range("A67").Value = "VALLEJO"
$ws.Range("A68").Value = "WILLOW CREEK"

# VALLEJO (row 67), like REDWOOD CITY above it, uses the wrap-text style.
$ws.Range("A67").WrapText = $true

# Fill in the latitude / longitude values for the inserted row.
$ws.Range("B56").Value = 34.966440261044198
$ws.Range("C56").Value = -120.657760310366

# Fill in the latitude / longitude values for the newly appended rows,
# using the same numeric display format ("0.0000") as the rest of the
# coordinate columns.
$ws.Range("B63").Value = 46.259540000000001
$ws.Range("C63").Value = -124.082083
$ws.Range("B64").Value = 40.727520494834501
$ws.Range("C64").Value = -124.22085681388
$ws.Range("B65").Value = 38.143129512729601
$ws.Range("C65").Value = -122.88032128675
$ws.Range("B66").Value = 42.738932436794897
$ws.Range("C66").Value = -124.49862438502601
$ws.Range("B67").Value = 38.089187955173301
$ws.Range("C67").Value = -122.29490950278699
$ws.Range("B68").Value = 36.622672000000001
$ws.Range("C68").Value = -121.885687

$ws.Range("B63:C68").NumberFormat = "0.0000"

# Match the final active selection from the source edit.
$ws.Range("F54").Select() | Out-Null
